$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "76.706.19"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.129.54"
$ws.Range("E3").Value = "  +5.74%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.17"
$ws.Range("E5").Value = "  +0.15%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "624.30"
$ws.Range("E6").Value = "  +4.45%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.218"
$ws.Range("E8").Value = "  +8.12%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  +0.77%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  +10.13%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.65%  "

# Row 12 - Toncoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +6.58%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.711.00"
$ws.Range("E13").Value = "  +5.98%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.36"
$ws.Range("E14").Value = "  +5.76%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +7.07%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "76.623.74"
$ws.Range("E16").Value = "  +0.30%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.137.13"
$ws.Range("E17").Value = "  +6.19%  "

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.76"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.24"
$ws.Range("E19").Value = "  +4.49%  "

# Row 20 - SuiNetwork
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.77"
$ws.Range("E20").Value = "  +20.84%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "405.26"
$ws.Range("E21").Value = "  +7.38%  "

# Row 22 - Polkadot
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.62"
$ws.Range("E22").Value = "  +6.83%  "

# Row 23 - LEO
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.54"
$ws.Range("E23").Value = "  +2.06%  "

# Row 24 - NEARProtocol -> WrappedeETH
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.296.68"
$ws.Range("E24").Value = "  +6.13%  "

# Row 25 - WrappedeETH -> NEARProtocol
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.64"
$ws.Range("E25").Value = "  +7.08%  "

# Row 26 - Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "74.03"
$ws.Range("E26").Value = "  +1.69%  "

# Row 27 - Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.49"
$ws.Range("E27").Value = "  +7.96%  "

# Row 28 - Dai
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.09%  "

# Row 29 - PEPE
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000112"
$ws.Range("E29").Value = "  +3.66%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  +0.21%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.54"
$ws.Range("E31").Value = "  -0.32%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +5.15%  "

# Row 33 - Bittensor
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "516.47"
$ws.Range("E33").Value = "  +3.64%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +7.08%  "

# Row 35 - Kaspa
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.134"
$ws.Range("E35").Value = "  +20.22%  "

# Row 36 - EthereumClassic
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "21.65"
$ws.Range("E36").Value = "  +5.99%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.07%  "

# Row 38 - Monero
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.77"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39 - Aave
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "195.79"
$ws.Range("E39").Value = "  +7.67%  "

# Row 40 - PolygonEcosystemToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.386"
$ws.Range("E40").Value = "  -1.63%  "

# Row 41 - WhiteBITCoin
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.07"
$ws.Range("E41").Value = "  +0.47%  "

# Row 42 - Cronos
$ws.Range("E42").Value = "  -4.59%  "

# Row 43 - RenderToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.33"
$ws.Range("E43").Value = "  +7.11%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  +0.04%  "

# Row 45 - Mantle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.806"
$ws.Range("E45").Value = "  +20.21%  "

# Row 46 - Stacks -> ImmutableX
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.30"
$ws.Range("E46").Value = "  +8.09%  "

# Row 47 - ImmutableX -> Stacks
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +4.11%  "

# Row 48 - OKB
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.19"
$ws.Range("E48").Value = "  +5.92%  "

# Row 49 - dogwifhat
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.48"
$ws.Range("E49").Value = "  +6.76%  "

# Row 50 - ARBITRUM
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.620"
$ws.Range("E50").Value = "  +4.75%  "

# Row 51 - Filecoin
$ws.Range("E51").Value = "  +3.58%  "
